# 47 BRAND FIX SUM
# Corrects the "Цена" (price, column Q) values on the "Загрузка" sheet
# for the product rows, and restores the post-edit selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> corrected price (column Q / 17) map, as captured from the fix.
$priceFixes = @{
    2=1551.55; 3=1451.45; 4=1601.6; 5=1851.85; 6=1851.85; 7=1551.55; 8=1501.5; 9=1451.45; 10=1401.4;
    11=1351.35; 12=1351.35; 13=1851.85; 14=1701.7; 15=1451.45; 16=1451.45; 17=1451.45; 18=1451.45;
    19=1351.35; 20=1351.35; 21=1351.35; 22=1351.35; 27=1701.7; 28=1701.7; 29=1601.6; 30=1551.55;
    31=1451.45; 32=1701.7; 33=1701.7; 34=1701.7; 35=1401.4; 36=1401.4; 37=1401.4; 38=1501.5; 39=1501.5;
    40=1501.5; 41=1501.5; 42=1501.5; 43=1501.5; 44=1501.5; 45=1501.5; 46=1501.5; 47=1501.5; 48=1501.5;
    49=1501.5; 50=1501.5; 51=1501.5; 52=1501.5; 53=1501.5; 54=1501.5; 55=1501.5; 56=1501.5; 57=1501.5;
    58=1501.5; 59=1451.45; 60=1451.45; 61=1451.45; 62=1451.45; 63=1451.45; 64=1451.45; 65=1451.45;
    66=1451.45; 67=1451.45; 68=1451.45; 69=1451.45; 70=1451.45; 71=1451.45; 72=1451.45; 73=1451.45;
    74=1451.45; 75=1801.8; 76=1551.55; 77=1601.6; 78=1351.35; 79=1351.35; 80=1351.35; 81=1351.35;
    82=1351.35; 83=1351.35; 84=1451.45; 85=1451.45; 86=1451.45; 87=1451.45; 88=1451.45; 89=1351.35;
    90=1351.35; 91=1351.35; 92=1351.35; 93=1351.35; 94=1501.5; 95=1501.5; 96=1501.5; 98=1351.35;
    99=1351.35; 101=1351.35; 102=1351.35; 103=1651.65; 104=1451.45; 105=1451.45; 106=1851.85;
    107=1851.85; 108=1851.85; 109=1351.35; 110=1351.35; 111=1351.35; 112=1351.35; 113=1351.35;
    114=1351.35; 115=1351.35; 116=1351.35; 118=1351.35; 119=1351.35; 120=1451.45; 121=1451.45;
    122=1451.45; 123=1351.35; 124=1351.35; 125=1501.5; 126=1501.5; 127=1501.5; 128=1501.5; 129=1501.5;
    130=1501.5; 131=1501.5; 132=1501.5; 133=1501.5
}

foreach ($row in $priceFixes.Keys) {
    $ws.Cells.Item($row, 17).Value = $priceFixes[$row]
}

# Restore the view/selection state left behind by the edit: column R
# (the full column) selected, with the frozen header row still in place.
$ws.Activate() | Out-Null
$ws.Range("R1:R1048576").Select() | Out-Null
